# Auto-generated edit script: update Typhon_Profits market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
#  LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the ALC/ARM/BSM/CRP/CUL/GSM/WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H29").Value = 1881.125
$ws.Range("J29").Value = 2135.5715
$ws.Range("L29").Value = 6406.7145
$ws.Range("N29").Value = -6968.7145

$ws.Range("H129").Value = 877.3958
$ws.Range("I129").Value = 665.6667
$ws.Range("J129").Value = 891.5111000000001
$ws.Range("K129").Value = 1997.0001
$ws.Range("L129").Value = 2674.5333
$ws.Range("M129").Value = 3002.9999
$ws.Range("N129").Value = -12674.5333

$ws.Range("H137").Value = 2218.4482
$ws.Range("I137").Value = 2559.6428
$ws.Range("J137").Value = 1900
$ws.Range("K137").Value = 7678.928400000001
$ws.Range("L137").Value = 5700
$ws.Range("M137").Value = -5128.928400000001
$ws.Range("N137").Value = -10800

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 2349.791
$ws.Range("I32").Value = 1438.3276
$ws.Range("J32").Value = 8223.666999999999
$ws.Range("K32").Value = 1438.3276
$ws.Range("L32").Value = 8223.666999999999
$ws.Range("M32").Value = -1151.3276
$ws.Range("N32").Value = -8797.666999999999

$ws.Range("H45").Value = 2643.353
$ws.Range("I45").Value = 2326.625
$ws.Range("J45").Value = 2924.889
$ws.Range("K45").Value = 2326.625
$ws.Range("L45").Value = 2924.889
$ws.Range("M45").Value = -1949.625
$ws.Range("N45").Value = -3678.889

$ws.Range("H61").Value = 4050.84
$ws.Range("I61").Value = 3809.2144
$ws.Range("K61").Value = 3809.2144
$ws.Range("M61").Value = -3597.2144

$ws.Range("H74").Value = 1045.125
$ws.Range("I74").Value = 675.125
$ws.Range("J74").Value = 1230.125
$ws.Range("K74").Value = 675.125
$ws.Range("L74").Value = 1230.125
$ws.Range("M74").Value = 198.875
$ws.Range("N74").Value = -2978.125

$ws.Range("H77").Value = 1045.125
$ws.Range("I77").Value = 675.125
$ws.Range("J77").Value = 1230.125
$ws.Range("K77").Value = 3375.625
$ws.Range("L77").Value = 6150.625
$ws.Range("M77").Value = 992.375
$ws.Range("N77").Value = -14886.625

$ws.Range("H122").Value = 2586.04
$ws.Range("I122").Value = 2393.6924
$ws.Range("J122").Value = 2794.4167
$ws.Range("K122").Value = 7181.0772
$ws.Range("L122").Value = 8383.250100000001
$ws.Range("M122").Value = -4731.0772
$ws.Range("N122").Value = -13283.2501

$ws.Range("H136").Value = 4050.84
$ws.Range("I136").Value = 3809.2144
$ws.Range("K136").Value = 11427.6432
$ws.Range("M136").Value = -8877.643199999999

# Row 70: leve fully crafted in-house (no market purchase) -> zero out price cols, drop profit cell(s)
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73: leve fully crafted in-house (no market purchase) -> zero out price cols, drop profit cell(s)
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H94").Value = 4292.4375
$ws.Range("I94").Value = 1789.8
$ws.Range("J94").Value = 5430
$ws.Range("K94").Value = 1789.8
$ws.Range("L94").Value = 5430
$ws.Range("M94").Value = -1338.8
$ws.Range("N94").Value = -6332

# Row 135: leve fully crafted in-house (no market purchase) -> zero out price cols, drop profit cell(s)
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H99").Value = 5172
$ws.Range("I99").Value = 4042.8572
$ws.Range("J99").Value = 6609.091
$ws.Range("K99").Value = 4042.8572
$ws.Range("L99").Value = 6609.091
$ws.Range("M99").Value = -2544.8572
$ws.Range("N99").Value = -9605.091

$ws.Range("H122").Value = 1331.6364
$ws.Range("J122").Value = 1216.6666
$ws.Range("L122").Value = 3649.9998
$ws.Range("N122").Value = -8549.9998

$ws.Range("H126").Value = 5172
$ws.Range("I126").Value = 4042.8572
$ws.Range("J126").Value = 6609.091
$ws.Range("K126").Value = 12128.5716
$ws.Range("L126").Value = 19827.273
$ws.Range("M126").Value = -9658.571599999999
$ws.Range("N126").Value = -24767.273

$ws.Range("H132").Value = 3920.6086
$ws.Range("I132").Value = 2279.2942
$ws.Range("J132").Value = 8571
$ws.Range("K132").Value = 6837.882599999999
$ws.Range("L132").Value = 25713
$ws.Range("M132").Value = -4307.882599999999
$ws.Range("N132").Value = -30773

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 650.7
$ws.Range("I5").Value = 497.6316
$ws.Range("J5").Value = 789.1905
$ws.Range("K5").Value = 1492.8948
$ws.Range("L5").Value = 2367.5715
$ws.Range("M5").Value = -1380.8948
$ws.Range("N5").Value = -2591.5715

$ws.Range("H23").Value = 697.1111
$ws.Range("I23").Value = 500.5
$ws.Range("J23").Value = 753.2857
$ws.Range("K23").Value = 1501.5
$ws.Range("L23").Value = 2259.8571
$ws.Range("M23").Value = -1266.5
$ws.Range("N23").Value = -2729.8571

$ws.Range("H68").Value = 1063.6666
$ws.Range("I68").Value = 530.6923
$ws.Range("J68").Value = 1410.1
$ws.Range("K68").Value = 1592.0769
$ws.Range("L68").Value = 4230.299999999999
$ws.Range("M68").Value = -781.0769
$ws.Range("N68").Value = -5852.299999999999

$ws.Range("H71").Value = 1063.6666
$ws.Range("I71").Value = 530.6923
$ws.Range("J71").Value = 1410.1
$ws.Range("K71").Value = 4776.2307
$ws.Range("L71").Value = 12690.9
$ws.Range("M71").Value = -720.2307000000001
$ws.Range("N71").Value = -20802.9

$ws.Range("H106").Value = 3480
$ws.Range("J106").Value = 3480
$ws.Range("L106").Value = 10440
$ws.Range("N106").Value = -12332

$ws.Range("H107").Value = 3750.75
$ws.Range("J107").Value = 720.9474
$ws.Range("L107").Value = 2162.8422
$ws.Range("N107").Value = -6002.8422

$ws.Range("H131").Value = 779.1
$ws.Range("J131").Value = 821.55914
$ws.Range("L131").Value = 2464.67742
$ws.Range("N131").Value = -12544.67742

$ws.Range("H135").Value = 650.7
$ws.Range("I135").Value = 497.6316
$ws.Range("J135").Value = 789.1905
$ws.Range("K135").Value = 4478.6844
$ws.Range("L135").Value = 7102.7145
$ws.Range("M135").Value = -1943.6844
$ws.Range("N135").Value = -12172.7145

# Row 33: leve fully crafted in-house (no market purchase) -> zero out price cols, drop profit cell(s)
$ws.Range("H33").Value = 650
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 650
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 3900
$ws.Range("N33").Value = -4466
$ws.Range("M33").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H102").Value = 2384.125
$ws.Range("I102").Value = 2381.8572
$ws.Range("K102").Value = 2381.8572
$ws.Range("M102").Value = -759.8571999999999

$ws.Range("H126").Value = 5510.0967
$ws.Range("I126").Value = 5649.9443
$ws.Range("K126").Value = 16949.8329
$ws.Range("M126").Value = -14479.8329

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872

$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360

$ws.Range("H132").Value = 1621.8485
$ws.Range("I132").Value = 1366.2307
$ws.Range("J132").Value = 2571.2856
$ws.Range("K132").Value = 4098.6921
$ws.Range("L132").Value = 7713.8568
$ws.Range("M132").Value = -1568.6921
$ws.Range("N132").Value = -12773.8568

$ws.Range("H136").Value = 1109.8667
$ws.Range("I136").Value = 637.05
$ws.Range("K136").Value = 1911.15
$ws.Range("M136").Value = 638.8500000000001

